$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1058.45
$ws.Range("J17").Value = 1058.45
$ws.Range("L17").Value = 3175.35
$ws.Range("N17").Value = -3511.35

$ws.Range("H64").Value = 4129.7617
$ws.Range("I64").Value = 4267.727
$ws.Range("J64").Value = 3978
$ws.Range("K64").Value = 4267.727
$ws.Range("L64").Value = 3978
$ws.Range("M64").Value = -4019.727
$ws.Range("N64").Value = -4474

$ws.Range("H67").Value = 4129.7617
$ws.Range("I67").Value = 4267.727
$ws.Range("J67").Value = 3978
$ws.Range("K67").Value = 4267.727
$ws.Range("L67").Value = 3978
$ws.Range("M67").Value = -3409.727
$ws.Range("N67").Value = -5694

$ws.Range("H96").Value = 808.6667
$ws.Range("I96").Value = 726
$ws.Range("K96").Value = 2178
$ws.Range("M96").Value = -805

$ws.Range("H132").Value = 13340307
$ws.Range("I132").Value = 19616722
$ws.Range("J132").Value = 2924.625
$ws.Range("K132").Value = 58850166
$ws.Range("L132").Value = 8773.875
$ws.Range("M132").Value = -58847636
$ws.Range("N132").Value = -13833.875

$ws.Range("H135").Value = 571.0345
$ws.Range("I135").Value = 317.03705
$ws.Range("K135").Value = 2853.33345
$ws.Range("M135").Value = -318.3334500000001

$ws.Range("H137").Value = 1312
$ws.Range("I137").Value = 850.8333
$ws.Range("J137").Value = 2142.1
$ws.Range("K137").Value = 2552.4999
$ws.Range("L137").Value = 6426.299999999999
$ws.Range("M137").Value = -2.499899999999798
$ws.Range("N137").Value = -11526.3

$ws.Range("H138").Value = 926771.0600000001
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 926771.0600000001
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 2780313.18
$ws.Range("N138").Value = -2790593.18
$ws.Range("M138").ClearContents()

$ws.Range("H141").Value = 593.2
$ws.Range("I141").Value = 593.2
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1779.6
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3400.4
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4166.3447
$ws.Range("I32").Value = 3971.7659
$ws.Range("J32").Value = 4997.727
$ws.Range("K32").Value = 3971.7659
$ws.Range("L32").Value = 4997.727
$ws.Range("M32").Value = -3684.7659
$ws.Range("N32").Value = -5571.727

$ws.Range("H61").Value = 1820.1177
$ws.Range("I61").Value = 1662.8
$ws.Range("K61").Value = 1662.8
$ws.Range("M61").Value = -1450.8

$ws.Range("H63").Value = 2050
$ws.Range("I63").Value = 2000
$ws.Range("K63").Value = 2000
$ws.Range("M63").Value = -1314

$ws.Range("H66").Value = 2050
$ws.Range("I66").Value = 2000
$ws.Range("K66").Value = 10000
$ws.Range("M66").Value = -6568

$ws.Range("H74").Value = 1595
$ws.Range("I74").Value = 735.1818
$ws.Range("J74").Value = 3171.3333
$ws.Range("K74").Value = 735.1818
$ws.Range("L74").Value = 3171.3333
$ws.Range("M74").Value = 138.8182
$ws.Range("N74").Value = -4919.3333

$ws.Range("H77").Value = 1595
$ws.Range("I77").Value = 735.1818
$ws.Range("J77").Value = 3171.3333
$ws.Range("K77").Value = 3675.909
$ws.Range("L77").Value = 15856.6665
$ws.Range("M77").Value = 692.0910000000003
$ws.Range("N77").Value = -24592.6665

$ws.Range("H122").Value = 2306
$ws.Range("I122").Value = 2306
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6918
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4468
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 1820.1177
$ws.Range("I136").Value = 1662.8
$ws.Range("K136").Value = 4988.4
$ws.Range("M136").Value = -2438.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4827.1333
$ws.Range("I86").Value = 5510.8
$ws.Range("J86").Value = 3459.8
$ws.Range("K86").Value = 5510.8
$ws.Range("L86").Value = 3459.8
$ws.Range("M86").Value = -4387.8
$ws.Range("N86").Value = -5705.8

$ws.Range("H89").Value = 4827.1333
$ws.Range("I89").Value = 5510.8
$ws.Range("J89").Value = 3459.8
$ws.Range("K89").Value = 27554
$ws.Range("L89").Value = 17299
$ws.Range("M89").Value = -21938
$ws.Range("N89").Value = -28531

$ws.Range("H107").Value = 1153.5769
$ws.Range("I107").Value = 968.35
$ws.Range("J107").Value = 1771
$ws.Range("K107").Value = 968.35
$ws.Range("L107").Value = 1771
$ws.Range("M107").Value = 951.65
$ws.Range("N107").Value = -5611

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1222.1177
$ws.Range("I31").Value = 853.2174
$ws.Range("J31").Value = 1993.4546
$ws.Range("K31").Value = 853.2174
$ws.Range("L31").Value = 1993.4546
$ws.Range("M31").Value = -558.2174
$ws.Range("N31").Value = -2583.4546

$ws.Range("H34").Value = 1222.1177
$ws.Range("I34").Value = 853.2174
$ws.Range("J34").Value = 1993.4546
$ws.Range("K34").Value = 853.2174
$ws.Range("L34").Value = 1993.4546
$ws.Range("M34").Value = -651.2174
$ws.Range("N34").Value = -2397.4546

$ws.Range("H62").Value = 10002320
$ws.Range("I62").Value = 2442.1052
$ws.Range("K62").Value = 2442.1052
$ws.Range("M62").Value = -1818.1052

$ws.Range("H65").Value = 10002320
$ws.Range("I65").Value = 2442.1052
$ws.Range("K65").Value = 12210.526
$ws.Range("M65").Value = -9090.526

$ws.Range("H99").Value = 1656.2142
$ws.Range("I99").Value = 1629.7693
$ws.Range("K99").Value = 1629.7693
$ws.Range("M99").Value = -131.7692999999999

$ws.Range("H107").Value = 527.65
$ws.Range("I107").Value = 424.2
$ws.Range("J107").Value = 838
$ws.Range("K107").Value = 424.2
$ws.Range("L107").Value = 838
$ws.Range("M107").Value = 1495.8
$ws.Range("N107").Value = -4678

$ws.Range("H109").Value = 22000.125
$ws.Range("J109").Value = 22000.125
$ws.Range("L109").Value = 22000.125
$ws.Range("N109").Value = -24080.125

$ws.Range("H126").Value = 1656.2142
$ws.Range("I126").Value = 1629.7693
$ws.Range("K126").Value = 4889.3079
$ws.Range("M126").Value = -2419.3079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1772.08
$ws.Range("I131").Value = 499.5
$ws.Range("J131").Value = 1798.051
$ws.Range("K131").Value = 1498.5
$ws.Range("L131").Value = 5394.153
$ws.Range("M131").Value = 3541.5
$ws.Range("N131").Value = -15474.153

$ws.Range("H139").Value = 1872.1316
$ws.Range("I139").Value = 2097.9
$ws.Range("J139").Value = 1621.2778
$ws.Range("K139").Value = 6293.700000000001
$ws.Range("L139").Value = 4863.8334
$ws.Range("M139").Value = -1153.700000000001
$ws.Range("N139").Value = -15143.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 5900
$ws.Range("I29").Value = 5900
$ws.Range("K29").Value = 5900
$ws.Range("M29").Value = -5610

$ws.Range("H46").Value = 12124.667
$ws.Range("I46").Value = 2599.6667
$ws.Range("J46").Value = 21649.666
$ws.Range("K46").Value = 2599.6667
$ws.Range("L46").Value = 21649.666
$ws.Range("M46").Value = -2443.6667
$ws.Range("N46").Value = -21961.666

$ws.Range("H70").Value = 20462816
$ws.Range("I70").Value = 14715579
$ws.Range("J70").Value = 40003420
$ws.Range("K70").Value = 14715579
$ws.Range("L70").Value = 40003420
$ws.Range("M70").Value = -14715309
$ws.Range("N70").Value = -40003960

$ws.Range("H73").Value = 20462816
$ws.Range("I73").Value = 14715579
$ws.Range("J73").Value = 40003420
$ws.Range("K73").Value = 14715579
$ws.Range("L73").Value = 40003420
$ws.Range("M73").Value = -14714643
$ws.Range("N73").Value = -40005292

$ws.Range("H80").Value = 3986.5715
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 4484.3335
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 4484.3335
$ws.Range("M80").Value = -2
$ws.Range("N80").Value = -6480.3335

$ws.Range("H83").Value = 3986.5715
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 4484.3335
$ws.Range("K83").Value = 5000
$ws.Range("L83").Value = 22421.6675
$ws.Range("M83").Value = -8
$ws.Range("N83").Value = -32405.6675

$ws.Range("H97").Value = 858.46155
$ws.Range("I97").Value = 846
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 846
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -350
$ws.Range("N97").Value = -1892

$ws.Range("H113").Value = 1483.6666
$ws.Range("I113").Value = 1397.7778
$ws.Range("J113").Value = 1612.5
$ws.Range("K113").Value = 1397.7778
$ws.Range("L113").Value = 1612.5
$ws.Range("M113").Value = 772.2221999999999
$ws.Range("N113").Value = -5952.5

$ws.Range("H122").Value = 1993.7894
$ws.Range("J122").Value = 2550
$ws.Range("L122").Value = 7650
$ws.Range("N122").Value = -12550

$ws.Range("H126").Value = 2236.25
$ws.Range("I126").Value = 1778
$ws.Range("K126").Value = 5334
$ws.Range("M126").Value = -2864

$ws.Range("H132").Value = 3012.6667
$ws.Range("I132").Value = 2982.5715
$ws.Range("J132").Value = 3027.7144
$ws.Range("K132").Value = 8947.7145
$ws.Range("L132").Value = 9083.143199999999
$ws.Range("M132").Value = -6417.7145
$ws.Range("N132").Value = -14143.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3744.5454
$ws.Range("I40").Value = 2748.5715
$ws.Range("K40").Value = 2748.5715
$ws.Range("M40").Value = -2612.5715

$ws.Range("H82").Value = 1972.9231
$ws.Range("J82").Value = 1949.5
$ws.Range("L82").Value = 1949.5
$ws.Range("N82").Value = -2671.5

$ws.Range("H85").Value = 1972.9231
$ws.Range("J85").Value = 1949.5
$ws.Range("L85").Value = 1949.5
$ws.Range("N85").Value = -4445.5

$ws.Range("H132").Value = 23656.543
$ws.Range("I132").Value = 1470.3478
$ws.Range("J132").Value = 45842.74
$ws.Range("K132").Value = 4411.0434
$ws.Range("L132").Value = 137528.22
$ws.Range("M132").Value = -1881.0434
$ws.Range("N132").Value = -142588.22

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 26001190
$ws.Range("I122").Value = 37144116
$ws.Range("J122").Value = 1033.3334
$ws.Range("K122").Value = 111432348
$ws.Range("L122").Value = 3100.0002
$ws.Range("M122").Value = -111429898
$ws.Range("N122").Value = -8000.0002

$ws.Range("H126").Value = 66668420
$ws.Range("J126").Value = 1900
$ws.Range("L126").Value = 5700
$ws.Range("N126").Value = -10640

$ws.Range("H132").Value = 4356.2383
$ws.Range("I132").Value = 5899.9
$ws.Range("J132").Value = 2952.9092
$ws.Range("K132").Value = 17699.7
$ws.Range("L132").Value = 8858.7276
$ws.Range("M132").Value = -15169.7
$ws.Range("N132").Value = -13918.7276
